$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text values are updated to match the refreshed crypto snapshot.
# D-column values that look like plain decimals need an apostrophe prefix
# (classic Excel "force text" input) so they stay text instead of being
# auto-parsed as numbers; the Style reset keeps the cell unstyled like the original.

$ws.Range("D2").Value = '43.083.82'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '2.306.11'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''301.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = '''98.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("D7").Value = '''0.524'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.07%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("D10").Value = '''35.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("D11").Value = '''0.0793'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '''17.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").Value = '2.663.94'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '2.305.51'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = '''0.790'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.38%  '
$ws.Range("D18").Value = '42.971.58'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").Value = '''13.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.72%  '
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("D22").Value = '''68.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("D23").Value = '''239.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = '''24.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").Value = '''167.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").Value = '''9.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("E30").Value = '  -13.34%  '
$ws.Range("D31").Value = '''33.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.98%  '
$ws.Range("D32").Value = '''5.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.39%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '''4.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("D35").Value = '''18.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.06%  '
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E40").Value = '  +2.02%  '
$ws.Range("E41").Value = '  -3.31%  '
$ws.Range("D42").Value = '2.004.96'
$ws.Range("E42").Value = '  +0.92%  '
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = '''10.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("D45").Value = '''2.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.23%  '
$ws.Range("D46").Value = '''17.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("E47").Value = '  -2.61%  '
$ws.Range("D48").Value = '''54.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("D49").Value = '2.529.49'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").Value = '''73.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.90%  '
